$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = $ws.ListObjects.Item("Table1")
$newRow = $table.ListRows.Add()

$ws.Range("A20").Value = "WEEK 2"
$ws.Range("B20").Value = "Rule based system"

$ws.Range("A20").HorizontalAlignment = $ws.Range("A19").HorizontalAlignment

$ws.Range("D20").Select()
